# Update the Kamlesh Nagarkoti match-log sheet: the existing single match
# row (row 2) is replaced with a new match, and five more match rows are
# appended below it (rows 3-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match data, in sheet (row) order.
# Columns: dateOfMatch, venueOfMatch, matchResult, ownTeam, opponentTeam,
#          playerName, runs, balls, numberOf4, numberOf6, sr
$matches = @(
    @(" Oct 12 2020", " Sharjah", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Kamlesh Nagarkoti ", "4", "7", "0", "0", "57.14"),
    @(" Oct 3 2020", " Sharjah", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", "Kamlesh Nagarkoti ", "3", "3", "0", "0", "100.00"),
    @(" Oct 7 2020", " Abu Dhabi", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", "Kamlesh Nagarkoti ", "0", "2", "0", "0", "0.00"),
    @(" Oct 26 2020", " Sharjah", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Kamlesh Nagarkoti ", "6", "13", "0", "0", "46.15"),
    @(" Nov 1 2020", " Dubai (DSC)", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Kamlesh Nagarkoti ", "1", "1", "0", "0", "100.00"),
    @(" Sep 30 2020", " Dubai (DSC)", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Kamlesh Nagarkoti ", "8", "5", "1", "0", "160.00")
)

$firstRow = 2
$lastRow = $firstRow + $matches.Length - 1

# The numeric-looking columns (runs/balls/4s/6s/sr) must stay text, just
# like the rest of the sheet, instead of being auto-coerced to numbers -
# format the block as Text before writing the values.
$numericCols = $ws.Range("G" + $firstRow + ":K" + $lastRow)
$numericCols.NumberFormat = "@"

for ($i = 0; $i -lt $matches.Length; $i++) {
    $row = $firstRow + $i
    $rowData = $matches[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

# Drop back to the sheet's normal style now that the text values are locked
# in, so we don't leave a stray "Text" number format behind on those cells.
$numericCols.Style = "Normal"
